$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the title cell (A1) with the corrected board name
$ws.Range("A1").Value = "Bill of Materials for 'Marote - M6-RF315 (Rev A)'"

# Correct the quantity-per-board values for rows 20 and 21 (were 0, should be 1)
$ws.Range("J20").Value = 1
$ws.Range("J21").Value = 1

# Reset the active selection to A2
$ws.Range("A2").Select()
